# Reorders the "Model" groups in rows 90-102 of the "Shop Visit Status" sheet:
# the HGT 1700 block (previously rows 100-102) moves to the top of this range
# (rows 90-92), and the APS1000 block (previously rows 90-91) moves to the
# bottom (rows 101-102). The APS5000 and GTCP331-500 blocks in between keep
# their relative order but shift up by 2 rows.
#
# Also updates the saved cell selection on the active sheet to D18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shop Visit Status")
if ($ws -eq $null) { $ws = $wb.ActiveSheet }

$srcRange = $ws.Range("A90:D102")
$src = $srcRange.Value2

$rows = $srcRange.Rows.Count
$cols = $srcRange.Columns.Count

# New source-row order (1-based, relative to the A90:D102 block) that
# reproduces the target layout:
#   HGT 1700 (was rows 100-102) -> rows 90-92
#   APS5000 + GTCP331-500 (was rows 92-99) -> rows 93-100 (unchanged order)
#   APS1000 (was rows 90-91) -> rows 101-102
$order = @(11, 12, 13, 3, 4, 5, 6, 7, 8, 9, 10, 1, 2)

$dst = New-Object 'object[,]' $rows, $cols
for ($r = 1; $r -le $rows; $r++) {
    $srcRowIndex = $order[$r - 1]
    for ($c = 1; $c -le $cols; $c++) {
        $dst[$r - 1, $c - 1] = $src[$srcRowIndex, $c]
    }
}

$srcRange.Value2 = $dst

# Restore the workbook's saved selection state (D10 -> D18 on the active
# sheet, which is the first/tab-selected worksheet).
$ws.Range("D18").Select()
